$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a flat weekly price-report table (row 1 = headers, A:T).
# The edit adds one more weekly batch of "Manzana" (apple) price quotes for
# "Terminal La Palmera de La Serena" dated 2021-09-09 (serial 44448).
# That batch is inserted as 9 new rows right before the existing row 958,
# pushing the previously-existing rows 958-984 down to 967-993.
# ---------------------------------------------------------------------------

# Insert 9 blank rows at 958..966 (pushes old 958-984 -> 967-993)
$ws.Range("A958:A966").EntireRow.Insert()

# Common (constant-across-the-batch) column values, matching the rest of
# this market/product block in the sheet.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$fecha     = 44448
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100104
$producto  = "Frutos de pepita"
$catId     = 100104002
$categoria = "Manzana"
$unidad    = "`$/bins (400 kilos)"
$kgUnidad  = 400

# Per-row variable data: Variedad, Calidad, Volumen, Precio min, Precio max,
# Precio promedio ponderado, Origen, Precio $/Kg
# NOTE: PowerShell hashtable keys are case-insensitive, so the row-number key
# is named "RowNum" (not "r"/"R") to avoid colliding with the Origen column.
$rows = @(
  @{ RowNum=958; K="Fuji royal";   L="Especial"; M=20; N=275000; O=280000; P=277500; R="Región de O'Higgins"; S=694 },
  @{ RowNum=959; K="Fuji royal";   L="Primera";  M=20; N=245000; O=250000; P=247500; R="Región de O'Higgins"; S=619 },
  @{ RowNum=960; K="Fuji royal";   L="Segunda";  M=16; N=215000; O=220000; P=217500; R="Región de O'Higgins"; S=544 },
  @{ RowNum=961; K="Granny Smith"; L="Especial"; M=20; N=255000; O=260000; P=257500; R="Región de O'Higgins"; S=644 },
  @{ RowNum=962; K="Granny Smith"; L="Primera";  M=20; N=225000; O=230000; P=227500; R="Región de O'Higgins"; S=569 },
  @{ RowNum=963; K="Granny Smith"; L="Segunda";  M=16; N=195000; O=200000; P=197500; R="Región de O'Higgins"; S=494 },
  @{ RowNum=964; K="Pink Lady";    L="Especial"; M=16; N=225000; O=230000; P=227500; R="Provincia de Curicó"; S=569 },
  @{ RowNum=965; K="Pink Lady";    L="Primera";  M=20; N=205000; O=210000; P=207500; R="Provincia de Curicó"; S=519 },
  @{ RowNum=966; K="Pink Lady";    L="Segunda";  M=16; N=180000; O=185000; P=182500; R="Provincia de Curicó"; S=456 }
)

foreach ($row in $rows) {
    $r = $row.RowNum
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
